$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23
$ws.Range("D23").Value = "[AP & mAP 내용 정리]`n분류기의 성능 평가를 위한 지난 포스팅(정밀도(Precision)와 재현율(Recall) 내용 정리)에 이어 이번"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2699"

# Row 28
$ws.Range("D28").Value = "개쩌는 Latex 캡쳐 도구 - Mathpix"
$ws.Range("E28").Value = "https://ropiens.tistory.com/84"

# Row 37
$ws.Range("D37").Value = "[Paper Review] Data-Distortion Guided Self-Distillation for Deep Neural Networks"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1437&mod=document&pageid=1"

# Row 40
$ws.Range("D40").Value = "가장 느린 현대 프로그래밍 언어"
$ws.Range("E40").Value = "https://insightcampus.co.kr/?kboard_content_redirect=12975"

# Row 46
$ws.Range("D46").Value = "증폭사지유도 (Augmented limb lead)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/374"
